# Generate Report for Archive
# - Update the "Ready for handoff" status text to "In Translation" across all sheets
# - Narrow the status/date columns (zh-cn / de-de columns) to their new report width

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update status text wherever it appears, on every worksheet ---
$sheetCount = $wb.Worksheets.Count
for ($s = 1; $s -le $sheetCount; $s++) {
    $ws = $wb.Worksheets.Item($s)
    $used = $ws.UsedRange
    $nrows = $used.Rows.Count
    $ncols = $used.Columns.Count
    for ($r = 1; $r -le $nrows; $r++) {
        for ($c = 1; $c -le $ncols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $cellText = [string]$cell.Value()
            if ($cellText -eq $oldStatus) {
                $cell.Value = $newStatus
            }
        }
    }
}

# --- Narrow the previously-wide status columns ---
$newColumnWidth = 12.5   # yields stored column width closest to the report's new narrower width

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth   # column E (zh-cn)
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth   # column F (de-de)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth       # column C (Status)
